$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 05:51"

# 2) Row 65 - Honduras: refresh counts
$ws.Cells.Item(65, 2).Value = 11258
$ws.Cells.Item(65, 3).Value = 519
$ws.Cells.Item(65, 4).Value = 1214
$ws.Cells.Item(65, 5).Value = 9695
$ws.Cells.Item(65, 7).Value = 6
$ws.Cells.Item(65, 8).Value = 349

# 3) Row 73 - Australia: refresh counts
$ws.Cells.Item(73, 2).Value = 7410
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(73, 4).Value = 6881
$ws.Cells.Item(73, 5).Value = 427

# 4) Row 164 - Mongolia: refresh counts
$ws.Cells.Item(164, 4).Value = 139
$ws.Cells.Item(164, 5).Value = 65

# 5) Rows 202/203 - swap Fiyi and Dominica (values are identical, just the
#    country names trade places)
$ws.Cells.Item(202, 1).Value = "Dominica"
$ws.Cells.Item(203, 1).Value = "Fiyi"

# 6) Rows 208/209 - swap Santa Sede and Islas Turcas y Caicos, including
#    their per-row data values
$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 1

$ws.Cells.Item(209, 1).Value = "Santa Sede"
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 8).Value = 0
